$d = $word.ActiveDocument

function Find-ParagraphIndexContaining($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($text)) {
            return $i
        }
    }
    return -1
}

# --- Change 1 -----------------------------------------------------------
# Insert a new non-bold paragraph ("Pour un développeur Backend comme on
# les appelle ") right after the "Pour un développeur Desktop ..."
# paragraph, leaving the existing empty bold paragraph that used to follow
# it intact (it is now pushed one slot further down).
$deskIdx = Find-ParagraphIndexContaining("Pour un développeur Desktop")
if ($deskIdx -eq -1) {
    throw "Could not locate the 'Pour un développeur Desktop' paragraph"
}
$deskPara = $d.Paragraphs.Item($deskIdx)
$deskPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($deskIdx + 1)
$insertion = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertion.InsertBefore("Pour un développeur Backend comme on les appelle ")

# --- Change 2 -----------------------------------------------------------
# Remove the empty bold paragraph that sits between the "Lens :" answer
# marker (after "Quelles sont les principales missions d'un développeur ?")
# and "Il vous est déjà arrivé de diriger ...".
$missionsIdx = Find-ParagraphIndexContaining("principales missions")
if ($missionsIdx -eq -1) {
    throw "Could not locate the 'principales missions' paragraph"
}
$lensIdx = $missionsIdx + 1
$emptyBoldIdx = $lensIdx + 1
$lensPara = $d.Paragraphs.Item($lensIdx)
$emptyBoldPara = $d.Paragraphs.Item($emptyBoldIdx)
if (-not ($lensPara.Range.Text.Contains("Lens"))) {
    throw "Unexpected paragraph layout around 'principales missions'"
}
$emptyBoldPara.Range.Delete()

# --- Change 3 -----------------------------------------------------------
# Merge the three footer runs that make up the address into a single run.
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$footerRange = $footer.Range
$footerRange.Find.Execute("253, Avenue Kimvula C/Bandalungwa, Kinshasa – République Démocratique du Congo", $true, $false, $false, $false, $false, $true, 1, $false, "253, Avenue Kimvula C/Bandalungwa, Kinshasa – République Démocratique du Congo", 2) | Out-Null

Write-Output "done"
